$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing GDP values (column B) for rows 12-30 with the revised
# FRED figures.
$ws.Range("B12").Value = 179376.59299999999
$ws.Range("B13").Value = 190252.19200000001
$ws.Range("B14").Value = 201324.05100000001
$ws.Range("B15").Value = 218381.45600000001
$ws.Range("B16").Value = 237504.83
$ws.Range("B17").Value = 251490.25399999999
$ws.Range("B18").Value = 260453.36300000001
$ws.Range("B19").Value = 253256.139
$ws.Range("B20").Value = 241248.18299999999
$ws.Range("B21").Value = 250416.44
$ws.Range("B22").Value = 254429.72099999999
$ws.Range("B23").Value = 261234.78899999999
$ws.Range("B24").Value = 273768.94699999999
$ws.Range("B25").Value = 288657.23
$ws.Range("B26").Value = 308286.88500000001
$ws.Range("B27").Value = 322562.83199999999
$ws.Range("B28").Value = 340554.337
$ws.Range("B29").Value = 360229.87
$ws.Range("B30").Value = 376603.20699999999

# Append the newly published 2020-01-01 observation as row 31, matching
# the date/number formatting already used for the rest of the series.
$ws.Range("A31").Value = 43831
$ws.Range("A31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B31").Value = 365051.489
$ws.Range("B31").NumberFormat = "0.000"

# Reflect the post-edit selection (whole A:B columns selected, cursor
# resting on A23).
$ws.Range("A1:B1048576").Select()
